# Auto-generated edit script applying scheduled market-data refresh values
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("ALC")
$ws.Range("H98").Value = 1173.9062
$ws.Range("I98").Value = 805.8929000000001
$ws.Range("K98").Value = 805.8929000000001
$ws.Range("M98").Value = 692.1070999999999

$ws.Range("H122").Value = 1173.9062
$ws.Range("I122").Value = 805.8929000000001
$ws.Range("K122").Value = 2417.6787
$ws.Range("M122").Value = 32.32129999999961

$ws = $wb.Worksheets("ARM")
$ws.Range("H8").Value = 23500
$ws.Range("I8").Value = 19000
$ws.Range("K8").Value = 19000
$ws.Range("M8").Value = -18856

$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").ClearContents()

$ws.Range("H13").Value = 60000
$ws.Range("J13").Value = 60000
$ws.Range("L13").Value = 60000
$ws.Range("N13").Value = -60288

$ws.Range("H61").Value = 2364.4
$ws.Range("I61").Value = 1396.2858
$ws.Range("J61").Value = 2885.6924
$ws.Range("K61").Value = 1396.2858
$ws.Range("L61").Value = 2885.6924
$ws.Range("M61").Value = -1184.2858
$ws.Range("N61").Value = -3309.6924

$ws.Range("H63").Value = 2004
$ws.Range("I63").Value = 2004
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 2004
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -1318
$ws.Range("N63").ClearContents()

$ws.Range("H66").Value = 2004
$ws.Range("I66").Value = 2004
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 10020
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -6588
$ws.Range("N66").ClearContents()

$ws.Range("H74").Value = 1365.9736
$ws.Range("I74").Value = 865.7586
$ws.Range("J74").Value = 2977.7778
$ws.Range("K74").Value = 865.7586
$ws.Range("L74").Value = 2977.7778
$ws.Range("M74").Value = 8.241399999999999
$ws.Range("N74").Value = -4725.7778

$ws.Range("H77").Value = 1365.9736
$ws.Range("I77").Value = 865.7586
$ws.Range("J77").Value = 2977.7778
$ws.Range("K77").Value = 4328.793
$ws.Range("L77").Value = 14888.889
$ws.Range("M77").Value = 39.20700000000033
$ws.Range("N77").Value = -23624.889

$ws.Range("H132").Value = 2819.1191
$ws.Range("I132").Value = 2639.3
$ws.Range("J132").Value = 3268.6667
$ws.Range("K132").Value = 7917.900000000001
$ws.Range("L132").Value = 9806.000100000001
$ws.Range("M132").Value = -5387.900000000001
$ws.Range("N132").Value = -14866.0001

$ws.Range("H136").Value = 2364.4
$ws.Range("I136").Value = 1396.2858
$ws.Range("J136").Value = 2885.6924
$ws.Range("K136").Value = 4188.857400000001
$ws.Range("L136").Value = 8657.0772
$ws.Range("M136").Value = -1638.857400000001
$ws.Range("N136").Value = -13757.0772

$ws = $wb.Worksheets("CRP")
$ws.Range("H31").Value = 1786.8334
$ws.Range("I31").Value = 1423.7858
$ws.Range("J31").Value = 3057.5
$ws.Range("K31").Value = 1423.7858
$ws.Range("L31").Value = 3057.5
$ws.Range("M31").Value = -1128.7858
$ws.Range("N31").Value = -3647.5

$ws.Range("H34").Value = 1786.8334
$ws.Range("I34").Value = 1423.7858
$ws.Range("J34").Value = 3057.5
$ws.Range("K34").Value = 1423.7858
$ws.Range("L34").Value = 3057.5
$ws.Range("M34").Value = -1221.7858
$ws.Range("N34").Value = -3461.5

$ws.Range("H58").Value = 1680.7858
$ws.Range("I58").Value = 1851.0588
$ws.Range("J58").Value = 1417.6364
$ws.Range("K58").Value = 1851.0588
$ws.Range("L58").Value = 1417.6364
$ws.Range("M58").Value = -1648.0588
$ws.Range("N58").Value = -1823.6364

$ws.Range("H132").Value = 1949.2858
$ws.Range("I132").Value = 1552.0555
$ws.Range("K132").Value = 4656.166499999999
$ws.Range("M132").Value = -2126.166499999999

$ws.Range("H136").Value = 1680.7858
$ws.Range("I136").Value = 1851.0588
$ws.Range("J136").Value = 1417.6364
$ws.Range("K136").Value = 5553.1764
$ws.Range("L136").Value = 4252.9092
$ws.Range("M136").Value = -3003.1764
$ws.Range("N136").Value = -9352.9092

$ws = $wb.Worksheets("CUL")
$ws.Range("H95").Value = 8950
$ws.Range("J95").Value = 8950
$ws.Range("L95").Value = 26850
$ws.Range("N95").Value = -30968

$ws.Range("H109").Value = 4275.7
$ws.Range("I109").Value = 1250.8
$ws.Range("J109").Value = 5284
$ws.Range("K109").Value = 3752.4
$ws.Range("L109").Value = 15852
$ws.Range("M109").Value = -2712.4
$ws.Range("N109").Value = -17932

$ws.Range("H122").Value = 1423.6774
$ws.Range("I122").Value = 1623.9
$ws.Range("K122").Value = 14615.1
$ws.Range("M122").Value = -12165.1

$ws = $wb.Worksheets("GSM")
$ws.Range("H28").Value = 15000
$ws.Range("J28").Value = 15000
$ws.Range("L28").Value = 15000
$ws.Range("N28").Value = -15384

$ws.Range("H102").Value = 2003.3103
$ws.Range("I102").Value = 1882.4348
$ws.Range("J102").Value = 2466.6667
$ws.Range("K102").Value = 1882.4348
$ws.Range("L102").Value = 2466.6667
$ws.Range("M102").Value = -260.4348
$ws.Range("N102").Value = -5710.6667

$ws.Range("H126").Value = 1501.7778
$ws.Range("I126").Value = 1370.4
$ws.Range("J126").Value = 1552.3077
$ws.Range("K126").Value = 4111.200000000001
$ws.Range("L126").Value = 4656.9231
$ws.Range("M126").Value = -1641.200000000001
$ws.Range("N126").Value = -9596.9231

$ws = $wb.Worksheets("LTW")
$ws.Range("H7").Value = 1530.7307
$ws.Range("I7").Value = 1145.9333
$ws.Range("J7").Value = 2055.4546
$ws.Range("K7").Value = 1145.9333
$ws.Range("L7").Value = 2055.4546
$ws.Range("M7").Value = -1033.9333
$ws.Range("N7").Value = -2279.4546

$ws.Range("H126").Value = 1530.7307
$ws.Range("I126").Value = 1145.9333
$ws.Range("J126").Value = 2055.4546
$ws.Range("K126").Value = 3437.7999
$ws.Range("L126").Value = 6166.3638
$ws.Range("M126").Value = -967.7999
$ws.Range("N126").Value = -11106.3638

$ws.Range("H132").Value = 7124.7666
$ws.Range("I132").Value = 3836.0527
$ws.Range("K132").Value = 11508.1581
$ws.Range("M132").Value = -8978.158100000001

$ws = $wb.Worksheets("WVR")
$ws.Range("H26").Value = 3990
$ws.Range("I26").Value = 2000
$ws.Range("J26").Value = 4487.5
$ws.Range("K26").Value = 2000
$ws.Range("L26").Value = 4487.5
$ws.Range("M26").Value = -1707
$ws.Range("N26").Value = -5073.5

$ws.Range("H29").Value = 23333.334
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()

$ws.Range("H132").Value = 1359.2354
$ws.Range("I132").Value = 1070.1333
$ws.Range("J132").Value = 1587.4736
$ws.Range("K132").Value = 3210.3999
$ws.Range("L132").Value = 4762.4208
$ws.Range("M132").Value = -680.3998999999999
$ws.Range("N132").Value = -9822.4208
